$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.5165012559121004, 0.7041422298085055, 0.6513016757365531, 0.8070326361037409, 0.6435116278483757, 14)
    3  = @(0.3104632072696188, 0.5383557835394773, 0.4470623602678305, 0.6686272207051029, 0.6163585566045678, 13)
    4  = @(0.3391134671465679, 0.5090805170705599, 0.4219980843258732, 0.6496137962865884, 0.578713438955796, 12)
    5  = @(0.3337917257302686, 0.4696055677140649, 0.2791590168366371, 0.5283550102314135, 0.4295536180009555, 11)
    6  = @(0.4045476033665492, 0.4608815459385528, 0.2775773177649309, 0.5268560693063438, 0.3557754509913868, 10)
    7  = @(0.3354180023350901, 0.401425754014993, 0.2011920589153808, 0.4485443778662049, 0.3158681298469871, 9)
    8  = @(0.3794834057612096, 0.444356733850047, 0.2440075877121284, 0.493971241786532, 0.3380615877347145, 8)
    9  = @(0.340609349648699, 0.4069218976430036, 0.2045529024993298, 0.452275250814512, 0.3213946727463862, 7)
    10 = @(0.3910026095914697, 0.4634125819798098, 0.2475590203832885, 0.49755303273449, 0.3370625692822314, 6)
    11 = @(0.3845419086351703, 0.4648443220891614, 0.2552201254894569, 0.505193156613841, 0.3663121039371297, 5)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("B$row").Value = $values[0]
    $ws.Range("C$row").Value = $values[1]
    $ws.Range("D$row").Value = $values[2]
    $ws.Range("E$row").Value = $values[3]
    $ws.Range("F$row").Value = $values[4]
    $ws.Range("G$row").Value = $values[5]
}
